# Auto-generated edit script: update F-column "想去人数" counts per diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 140
$ws.Range("F6").Value = 285
$ws.Range("F7").Value = 13267
$ws.Range("F9").Value = 124
$ws.Range("F10").Value = 312
$ws.Range("F11").Value = 5097
$ws.Range("F13").Value = 3614
$ws.Range("F14").Value = 52
$ws.Range("F18").Value = 128
$ws.Range("F23").Value = 89
$ws.Range("F24").Value = 113
$ws.Range("F25").Value = 4591
$ws.Range("F27").Value = 1967
$ws.Range("F29").Value = 278
$ws.Range("F30").Value = 7165
$ws.Range("F33").Value = 2148
$ws.Range("F34").Value = 2071
$ws.Range("F35").Value = 1312
$ws.Range("F36").Value = 124
$ws.Range("F37").Value = 1117
$ws.Range("F39").Value = 7
$ws.Range("F40").Value = 236
$ws.Range("F41").Value = 229
$ws.Range("F43").Value = 11
$ws.Range("F45").Value = 1255
$ws.Range("F46").Value = 1886
$ws.Range("F48").Value = 177
$ws.Range("F49").Value = 1190

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F13").Value = 928

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 494
$ws.Range("F3").Value = 669
$ws.Range("F4").Value = 48

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 140
$ws.Range("F5").Value = 494
$ws.Range("F6").Value = 669
$ws.Range("F7").Value = 285
$ws.Range("F8").Value = 13267
$ws.Range("F10").Value = 312
$ws.Range("F11").Value = 5097
$ws.Range("F12").Value = 3614
$ws.Range("F13").Value = 52
$ws.Range("F16").Value = 128
$ws.Range("F21").Value = 89
$ws.Range("F23").Value = 113
$ws.Range("F24").Value = 4592
$ws.Range("F26").Value = 1967
$ws.Range("F28").Value = 278
$ws.Range("F29").Value = 7166
$ws.Range("F33").Value = 2148
$ws.Range("F34").Value = 2071
$ws.Range("F35").Value = 1312
$ws.Range("F36").Value = 124
$ws.Range("F37").Value = 1117
$ws.Range("F38").Value = 7
$ws.Range("F39").Value = 236
$ws.Range("F40").Value = 229
$ws.Range("F44").Value = 1255
$ws.Range("F45").Value = 1886
$ws.Range("F48").Value = 177
$ws.Range("F49").Value = 1190
